$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row A1:F1 gets a wrap-text cell style (new cellXfs entry).
$ws.Range("A1:F1").WrapText = $true

# Fix the Sensor 1 reading at row 7 (ship class / sensed-state test data).
$ws.Range("B7").Value = 1

# F2:F12 formulas now compare Sensed State vs Truth State directly
# instead of gating on the "sensor working" flag.
$ws.Range("F2").Formula = "=IF(B2=D2, 1, 0)"
$ws.Range("F3").Formula = "=IF(B3=D3, 1, 0)"
$ws.Range("F4").Formula = "=IF(B4=D4, 1, 0)"
$ws.Range("F5").Formula = "=IF(B5=D5, 1, 0)"
$ws.Range("F6").Formula = "=IF(B6=D6, 1, 0)"
$ws.Range("F7").Formula = "=IF(B7=D7, 1, 0)"
$ws.Range("F8").Formula = "=IF(B8=D8, 1, 0)"
$ws.Range("F9").Formula = "=IF(B9=D9, 1, 0)"
$ws.Range("F10").Formula = "=IF(B10=D10, 1, 0)"
$ws.Range("F11").Formula = "=IF(B11=D11, 1, 0)"
$ws.Range("F12").Formula = "=IF(B12=D12, 1, 0)"

# Conditional formatting: red->green color scale on the sensor-working
# column (E) and the match/no-match output column (F).
$rngE = $ws.Range("E2:E12")
$null = $rngE.FormatConditions.AddColorScale(2)
$cfE = $rngE.FormatConditions(1)
$cfE.ColorScaleCriteria(1).Type = 1
$cfE.ColorScaleCriteria(1).Value = 0
$cfE.ColorScaleCriteria(1).FormatColor.Color = 253
$cfE.ColorScaleCriteria(2).Type = 2
$cfE.ColorScaleCriteria(2).Value = 0
$cfE.ColorScaleCriteria(2).FormatColor.Color = 64768

$rngF = $ws.Range("F2:F12")
$null = $rngF.FormatConditions.AddColorScale(2)
$cfF = $rngF.FormatConditions(1)
$cfF.ColorScaleCriteria(1).Type = 1
$cfF.ColorScaleCriteria(1).Value = 0
$cfF.ColorScaleCriteria(1).FormatColor.Color = 253
$cfF.ColorScaleCriteria(2).Type = 2
$cfF.ColorScaleCriteria(2).Value = 0
$cfF.ColorScaleCriteria(2).FormatColor.Color = 64768
